$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 141
$ws.Range("I9").Value = 160.2
$ws.Range("K9").Value = 160.2
$ws.Range("M9").Value = 8.800000000000011
$ws.Range("H28").Value = 1544.9584
$ws.Range("I28").Value = 1498.25
$ws.Range("J28").Value = 1638.375
$ws.Range("K28").Value = 1498.25
$ws.Range("L28").Value = 1638.375
$ws.Range("M28").Value = -1013.25
$ws.Range("N28").Value = -2608.375
$ws.Range("H31").Value = 40864.4
$ws.Range("I31").Value = 40864.4
$ws.Range("K31").Value = 122593.2
$ws.Range("M31").Value = -122363.2
$ws.Range("H40").Value = 8395.4
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 8395.4
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 8395.4
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -8745.4
$ws.Range("H41").Value = 993.8
$ws.Range("I41").Value = 718.6667
$ws.Range("J41").Value = 1218.909
$ws.Range("K41").Value = 718.6667
$ws.Range("L41").Value = 1218.909
$ws.Range("M41").Value = -278.6667
$ws.Range("N41").Value = -2098.909
$ws.Range("H70").Value = 4127.5713
$ws.Range("I70").Value = 696
$ws.Range("J70").Value = 4299.15
$ws.Range("K70").Value = 2088
$ws.Range("L70").Value = 12897.45
$ws.Range("M70").Value = -1818
$ws.Range("N70").Value = -13437.45
$ws.Range("H73").Value = 4127.5713
$ws.Range("I73").Value = 696
$ws.Range("J73").Value = 4299.15
$ws.Range("K73").Value = 2088
$ws.Range("L73").Value = 12897.45
$ws.Range("M73").Value = -1152
$ws.Range("N73").Value = -14769.45
$ws.Range("H74").Value = 9423.076999999999
$ws.Range("I74").Value = 6900
$ws.Range("K74").Value = 6900
$ws.Range("M74").Value = -5964
$ws.Range("H77").Value = 9423.076999999999
$ws.Range("I77").Value = 6900
$ws.Range("K77").Value = 34500
$ws.Range("M77").Value = -29820
$ws.Range("H82").Value = 6567.375
$ws.Range("I82").Value = 4648.4287
$ws.Range("K82").Value = 13945.2861
$ws.Range("M82").Value = -13539.2861
$ws.Range("H85").Value = 6567.375
$ws.Range("I85").Value = 4648.4287
$ws.Range("K85").Value = 13945.2861
$ws.Range("M85").Value = -12541.2861
$ws.Range("H99").Value = 77565.62
$ws.Range("I99").Value = 265.4
$ws.Range("J99").Value = 125878.25
$ws.Range("K99").Value = 796.1999999999999
$ws.Range("L99").Value = 377634.75
$ws.Range("M99").Value = 701.8000000000001
$ws.Range("N99").Value = -380630.75
$ws.Range("H100").Value = 4990.364
$ws.Range("I100").Value = 2799.25
$ws.Range("J100").Value = 10833.333
$ws.Range("K100").Value = 2799.25
$ws.Range("L100").Value = 10833.333
$ws.Range("M100").Value = -2258.25
$ws.Range("N100").Value = -11915.333
$ws.Range("H101").Value = 2172.1428
$ws.Range("I101").Value = 1196
$ws.Range("J101").Value = 3473.6667
$ws.Range("K101").Value = 3588
$ws.Range("L101").Value = 10421.0001
$ws.Range("M101").Value = -1966
$ws.Range("N101").Value = -13665.0001
$ws.Range("H106").Value = 4343.375
$ws.Range("I106").Value = 4372.933
$ws.Range("K106").Value = 4372.933
$ws.Range("M106").Value = -3741.933
$ws.Range("H113").Value = 5647.6
$ws.Range("I113").Value = 3918.4
$ws.Range("J113").Value = 6800.4
$ws.Range("K113").Value = 3918.4
$ws.Range("L113").Value = 6800.4
$ws.Range("M113").Value = -664.4000000000001
$ws.Range("N113").Value = -13308.4
$ws.Range("H125").Value = 11279.9
$ws.Range("I125").Value = 1176.1428
$ws.Range("J125").Value = 34855.332
$ws.Range("K125").Value = 10585.2852
$ws.Range("L125").Value = 313697.988
$ws.Range("M125").Value = -8125.2852
$ws.Range("N125").Value = -318617.988
$ws.Range("H126").Value = 69937.5
$ws.Range("J126").Value = 69937.5
$ws.Range("L126").Value = 69937.5
$ws.Range("N126").Value = -79817.5
$ws.Range("H127").Value = 6956.8335
$ws.Range("I127").Value = 1168.2
$ws.Range("J127").Value = 35900
$ws.Range("K127").Value = 3504.6
$ws.Range("L127").Value = 107700
$ws.Range("M127").Value = 1455.4
$ws.Range("N127").Value = -117620
$ws.Range("H129").Value = 2971232
$ws.Range("I129").Value = 2971232
$ws.Range("K129").Value = 8913696
$ws.Range("M129").Value = -8908696
$ws.Range("H130").Value = 100000
$ws.Range("J130").Value = 100000
$ws.Range("L130").Value = 100000
$ws.Range("N130").Value = -110040
$ws.Range("H131").Value = 8142.125
$ws.Range("I131").Value = 5689.8335
$ws.Range("K131").Value = 17069.5005
$ws.Range("M131").Value = -12029.5005
$ws.Range("H132").Value = 36504.2
$ws.Range("I132").Value = 2612.5
$ws.Range("K132").Value = 7837.5
$ws.Range("M132").Value = -5307.5
$ws.Range("H133").Value = 70000
$ws.Range("J133").Value = 70000
$ws.Range("L133").Value = 70000
$ws.Range("N133").Value = -80120
$ws.Range("H135").Value = 2103
$ws.Range("I135").Value = 2151.7273
$ws.Range("K135").Value = 19365.5457
$ws.Range("M135").Value = -16830.5457
$ws.Range("H137").Value = 1633.258
$ws.Range("I137").Value = 1151.5238
$ws.Range("K137").Value = 3454.5714
$ws.Range("M137").Value = -904.5713999999998
$ws.Range("H138").Value = 1903.2838
$ws.Range("I138").Value = 1004.871
$ws.Range("J138").Value = 2550.9768
$ws.Range("K138").Value = 3014.613
$ws.Range("L138").Value = 7652.930399999999
$ws.Range("M138").Value = 2125.387
$ws.Range("N138").Value = -17932.9304

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 12194.286
$ws.Range("I2").Value = 15642.904
$ws.Range("K2").Value = 15642.904
$ws.Range("M2").Value = -15529.904
$ws.Range("H17").Value = 200
$ws.Range("I17").Value = 200
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 200
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -27
$ws.Range("N17").ClearContents()
$ws.Range("H32").Value = 8353.823
$ws.Range("I32").Value = 6865
$ws.Range("J32").Value = 16989
$ws.Range("K32").Value = 6865
$ws.Range("L32").Value = 16989
$ws.Range("M32").Value = -6578
$ws.Range("N32").Value = -17563
$ws.Range("H45").Value = 2115.4707
$ws.Range("I45").Value = 1377.2
$ws.Range("J45").Value = 3170.1428
$ws.Range("K45").Value = 1377.2
$ws.Range("L45").Value = 3170.1428
$ws.Range("M45").Value = -1000.2
$ws.Range("N45").Value = -3924.1428
$ws.Range("H61").Value = 3825.842
$ws.Range("I61").Value = 1545.7778
$ws.Range("K61").Value = 1545.7778
$ws.Range("M61").Value = -1333.7778
$ws.Range("H63").Value = 5123.8887
$ws.Range("I63").Value = 7323.2
$ws.Range("K63").Value = 7323.2
$ws.Range("M63").Value = -6637.2
$ws.Range("H66").Value = 5123.8887
$ws.Range("I66").Value = 7323.2
$ws.Range("K66").Value = 36616
$ws.Range("M66").Value = -33184
$ws.Range("H74").Value = 1913.6333
$ws.Range("I74").Value = 1940.64
$ws.Range("K74").Value = 1940.64
$ws.Range("M74").Value = -1066.64
$ws.Range("H77").Value = 1913.6333
$ws.Range("I77").Value = 1940.64
$ws.Range("K77").Value = 9703.200000000001
$ws.Range("M77").Value = -5335.200000000001
$ws.Range("H110").Value = 2581.9
$ws.Range("I110").Value = 2581.9
$ws.Range("K110").Value = 2581.9
$ws.Range("M110").Value = -536.9000000000001
$ws.Range("H116").Value = 12194.286
$ws.Range("I116").Value = 15642.904
$ws.Range("K116").Value = 15642.904
$ws.Range("M116").Value = -13348.904
$ws.Range("H132").Value = 2814.3
$ws.Range("I132").Value = 2103.95
$ws.Range("K132").Value = 6311.849999999999
$ws.Range("M132").Value = -3781.849999999999
$ws.Range("H136").Value = 3825.842
$ws.Range("I136").Value = 1545.7778
$ws.Range("K136").Value = 4637.3334
$ws.Range("M136").Value = -2087.3334

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 126374.25
$ws.Range("J4").Value = 1000000
$ws.Range("L4").Value = 1000000
$ws.Range("N4").Value = -1000224
$ws.Range("H31").Value = 4261.778
$ws.Range("I31").Value = 1440.25
$ws.Range("J31").Value = 5067.9287
$ws.Range("K31").Value = 1440.25
$ws.Range("L31").Value = 5067.9287
$ws.Range("M31").Value = -1145.25
$ws.Range("N31").Value = -5657.9287
$ws.Range("H34").Value = 4261.778
$ws.Range("I34").Value = 1440.25
$ws.Range("J34").Value = 5067.9287
$ws.Range("K34").Value = 1440.25
$ws.Range("L34").Value = 5067.9287
$ws.Range("M34").Value = -1238.25
$ws.Range("N34").Value = -5471.9287
$ws.Range("H58").Value = 2831.0527
$ws.Range("I58").Value = 2491.2
$ws.Range("J58").Value = 3208.6667
$ws.Range("K58").Value = 2491.2
$ws.Range("L58").Value = 3208.6667
$ws.Range("M58").Value = -2288.2
$ws.Range("N58").Value = -3614.6667
$ws.Range("H62").Value = 6729.5
$ws.Range("I62").Value = 7110.5557
$ws.Range("J62").Value = 3300
$ws.Range("K62").Value = 7110.5557
$ws.Range("L62").Value = 3300
$ws.Range("M62").Value = -6486.5557
$ws.Range("N62").Value = -4548
$ws.Range("H65").Value = 6729.5
$ws.Range("I65").Value = 7110.5557
$ws.Range("J65").Value = 3300
$ws.Range("K65").Value = 35552.7785
$ws.Range("L65").Value = 16500
$ws.Range("M65").Value = -32432.7785
$ws.Range("N65").Value = -22740
$ws.Range("H86").Value = 103326.336
$ws.Range("I86").Value = 103326.336
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 103326.336
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -102203.336
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 103326.336
$ws.Range("I89").Value = 103326.336
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 516631.68
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -511015.68
$ws.Range("N89").ClearContents()
$ws.Range("H107").Value = 6477.6665
$ws.Range("I107").Value = 939.38464
$ws.Range("J107").Value = 20877.2
$ws.Range("K107").Value = 939.38464
$ws.Range("L107").Value = 20877.2
$ws.Range("M107").Value = 980.61536
$ws.Range("N107").Value = -24717.2
$ws.Range("H122").Value = 382500.38
$ws.Range("I122").Value = 730758.9399999999
$ws.Range("K122").Value = 2192276.82
$ws.Range("M122").Value = -2189826.82
$ws.Range("H134").Value = 2951.0806
$ws.Range("I134").Value = 1776.0233
$ws.Range("K134").Value = 5328.0699
$ws.Range("M134").Value = -2793.0699
$ws.Range("H136").Value = 2831.0527
$ws.Range("I136").Value = 2491.2
$ws.Range("J136").Value = 3208.6667
$ws.Range("K136").Value = 7473.599999999999
$ws.Range("L136").Value = 9626.000100000001
$ws.Range("M136").Value = -4923.599999999999
$ws.Range("N136").Value = -14726.0001

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 73.125
$ws.Range("I2").Value = 13.25
$ws.Range("K2").Value = 79.5
$ws.Range("M2").Value = 33.5
$ws.Range("H5").Value = 8044.2856
$ws.Range("I5").Value = 4044.2354
$ws.Range("J5").Value = 14226.182
$ws.Range("K5").Value = 12132.7062
$ws.Range("L5").Value = 42678.546
$ws.Range("M5").Value = -12020.7062
$ws.Range("N5").Value = -42902.546
$ws.Range("H26").Value = 133.66667
$ws.Range("I26").Value = 125.375
$ws.Range("J26").Value = 200
$ws.Range("K26").Value = 376.125
$ws.Range("L26").Value = 600
$ws.Range("M26").Value = -88.125
$ws.Range("N26").Value = -1176
$ws.Range("H31").Value = 420
$ws.Range("I31").Value = 420
$ws.Range("K31").Value = 1260
$ws.Range("M31").Value = -972
$ws.Range("H60").Value = 734.8
$ws.Range("I60").Value = 451.66666
$ws.Range("J60").Value = 824.2105
$ws.Range("K60").Value = 1354.99998
$ws.Range("L60").Value = 2472.6315
$ws.Range("M60").Value = -1103.99998
$ws.Range("N60").Value = -2974.6315
$ws.Range("H104").Value = 4923.375
$ws.Range("J104").Value = 5981.1665
$ws.Range("L104").Value = 17943.4995
$ws.Range("N104").Value = -23185.4995
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()
$ws.Range("H109").Value = 2007.3334
$ws.Range("I109").Value = 1511
$ws.Range("K109").Value = 4533
$ws.Range("M109").Value = -3493
$ws.Range("H116").Value = 2421.0908
$ws.Range("I116").Value = 1380.3334
$ws.Range("J116").Value = 3670
$ws.Range("K116").Value = 4141.0002
$ws.Range("L116").Value = 11010
$ws.Range("M116").Value = -699.0002000000004
$ws.Range("N116").Value = -17894
$ws.Range("H131").Value = 5265.4707
$ws.Range("I131").Value = 3116.5557
$ws.Range("J131").Value = 7683
$ws.Range("K131").Value = 9349.667099999999
$ws.Range("L131").Value = 23049
$ws.Range("M131").Value = -4309.667099999999
$ws.Range("N131").Value = -33129
$ws.Range("H134").Value = 5293.5
$ws.Range("I134").Value = 1278.7778
$ws.Range("J134").Value = 12520
$ws.Range("K134").Value = 3836.3334
$ws.Range("L134").Value = 37560
$ws.Range("M134").Value = 1233.6666
$ws.Range("N134").Value = -47700
$ws.Range("H135").Value = 8044.2856
$ws.Range("I135").Value = 4044.2354
$ws.Range("J135").Value = 14226.182
$ws.Range("K135").Value = 36398.1186
$ws.Range("L135").Value = 128035.638
$ws.Range("M135").Value = -33863.1186
$ws.Range("N135").Value = -133105.638
$ws.Range("H140").Value = 1997.4
$ws.Range("I140").Value = 1997.4
$ws.Range("K140").Value = 5992.200000000001
$ws.Range("M140").Value = -812.2000000000007

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 313.30768
$ws.Range("I2").Value = 315.35715
$ws.Range("J2").Value = 310.91666
$ws.Range("K2").Value = 315.35715
$ws.Range("L2").Value = 310.91666
$ws.Range("M2").Value = -202.35715
$ws.Range("N2").Value = -536.91666
$ws.Range("H7").Value = 3632332.5
$ws.Range("J7").Value = 3632332.5
$ws.Range("L7").Value = 3632332.5
$ws.Range("N7").Value = -3632556.5
$ws.Range("H8").Value = 3632332.5
$ws.Range("J8").Value = 3632332.5
$ws.Range("L8").Value = 3632332.5
$ws.Range("N8").Value = -3632610.5
$ws.Range("H10").Value = 5935.6
$ws.Range("J10").Value = 4900
$ws.Range("L10").Value = 4900
$ws.Range("N10").Value = -5238
$ws.Range("H18").Value = 5000
$ws.Range("J18").Value = 5000
$ws.Range("L18").Value = 5000
$ws.Range("N18").Value = -5586
$ws.Range("H70").Value = 383703.34
$ws.Range("I70").Value = 565555
$ws.Range("J70").Value = 20000
$ws.Range("K70").Value = 565555
$ws.Range("L70").Value = 20000
$ws.Range("M70").Value = -565285
$ws.Range("N70").Value = -20540
$ws.Range("H73").Value = 383703.34
$ws.Range("I73").Value = 565555
$ws.Range("J73").Value = 20000
$ws.Range("K73").Value = 565555
$ws.Range("L73").Value = 20000
$ws.Range("M73").Value = -564619
$ws.Range("N73").Value = -21872
$ws.Range("H80").Value = 126967.78
$ws.Range("I80").Value = 224762.2
$ws.Range("K80").Value = 224762.2
$ws.Range("M80").Value = -223764.2
$ws.Range("H83").Value = 126967.78
$ws.Range("I83").Value = 224762.2
$ws.Range("K83").Value = 1123811
$ws.Range("M83").Value = -1118819
$ws.Range("H97").Value = 678.3333
$ws.Range("I97").Value = 703.7143
$ws.Range("J97").Value = 589.5
$ws.Range("K97").Value = 703.7143
$ws.Range("L97").Value = 589.5
$ws.Range("M97").Value = -207.7143
$ws.Range("N97").Value = -1581.5
$ws.Range("H102").Value = 1290.8182
$ws.Range("I102").Value = 1028.8
$ws.Range("K102").Value = 1028.8
$ws.Range("M102").Value = 593.2
$ws.Range("H113").Value = 9791
$ws.Range("I113").Value = 6421.2
$ws.Range("K113").Value = 6421.2
$ws.Range("M113").Value = -4251.2
$ws.Range("H122").Value = 5202
$ws.Range("I122").Value = 2258.25
$ws.Range("J122").Value = 6183.25
$ws.Range("K122").Value = 6774.75
$ws.Range("L122").Value = 18549.75
$ws.Range("M122").Value = -4324.75
$ws.Range("N122").Value = -23449.75
$ws.Range("H126").Value = 5750
$ws.Range("I126").Value = 3500
$ws.Range("J126").Value = 6500
$ws.Range("K126").Value = 10500
$ws.Range("L126").Value = 19500
$ws.Range("M126").Value = -8030
$ws.Range("N126").Value = -24440
$ws.Range("H132").Value = 2373.1843
$ws.Range("I132").Value = 1266.138
$ws.Range("K132").Value = 3798.414
$ws.Range("M132").Value = -1268.414
$ws.Range("H136").Value = 11202.9
$ws.Range("J136").Value = 11202.9
$ws.Range("L136").Value = 33608.7
$ws.Range("N136").Value = -38708.7

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3579.1155
$ws.Range("I7").Value = 2294.2354
$ws.Range("J7").Value = 6006.1113
$ws.Range("K7").Value = 2294.2354
$ws.Range("L7").Value = 6006.1113
$ws.Range("M7").Value = -2182.2354
$ws.Range("N7").Value = -6230.1113
$ws.Range("H20").Value = 97997.39999999999
$ws.Range("J20").Value = 97997.39999999999
$ws.Range("L20").Value = 97997.39999999999
$ws.Range("N20").Value = -98449.39999999999
$ws.Range("H22").Value = 690.8823
$ws.Range("J22").Value = 997.5
$ws.Range("L22").Value = 997.5
$ws.Range("N22").Value = -1587.5
$ws.Range("H26").Value = 17000
$ws.Range("J26").Value = 17000
$ws.Range("L26").Value = 17000
$ws.Range("N26").Value = -17590
$ws.Range("H27").Value = 690.8823
$ws.Range("J27").Value = 997.5
$ws.Range("L27").Value = 997.5
$ws.Range("N27").Value = -1211.5
$ws.Range("H40").Value = 11691.294
$ws.Range("I40").Value = 16621
$ws.Range("J40").Value = 8240.5
$ws.Range("K40").Value = 16621
$ws.Range("L40").Value = 8240.5
$ws.Range("M40").Value = -16485
$ws.Range("N40").Value = -8512.5
$ws.Range("H46").Value = 3132.6667
$ws.Range("I46").Value = 2199
$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 2199
$ws.Range("L46").Value = 5000
$ws.Range("M46").Value = -2011
$ws.Range("N46").Value = -5376
$ws.Range("H55").Value = 1088.4348
$ws.Range("I55").Value = 1254.4706
$ws.Range("J55").Value = 618
$ws.Range("K55").Value = 1254.4706
$ws.Range("L55").Value = 618
$ws.Range("M55").Value = -1081.4706
$ws.Range("N55").Value = -964
$ws.Range("H61").Value = 2498.6047
$ws.Range("I61").Value = 2106.4194
$ws.Range("K61").Value = 2106.4194
$ws.Range("M61").Value = -1904.4194
$ws.Range("H68").Value = 5564.5
$ws.Range("I68").Value = 5218.625
$ws.Range("J68").Value = 5795.0835
$ws.Range("K68").Value = 5218.625
$ws.Range("L68").Value = 5795.0835
$ws.Range("M68").Value = -4469.625
$ws.Range("N68").Value = -7293.0835
$ws.Range("H71").Value = 5564.5
$ws.Range("I71").Value = 5218.625
$ws.Range("J71").Value = 5795.0835
$ws.Range("K71").Value = 26093.125
$ws.Range("L71").Value = 28975.4175
$ws.Range("M71").Value = -22349.125
$ws.Range("N71").Value = -36463.4175
$ws.Range("H102").Value = 8373.75
$ws.Range("J102").Value = 8373.75
$ws.Range("L102").Value = 8373.75
$ws.Range("N102").Value = -14863.75
$ws.Range("H113").Value = 2498.6047
$ws.Range("I113").Value = 2106.4194
$ws.Range("K113").Value = 2106.4194
$ws.Range("M113").Value = 63.58059999999978
$ws.Range("H122").Value = 4472.9375
$ws.Range("I122").Value = 2558.8572
$ws.Range("J122").Value = 5961.6665
$ws.Range("K122").Value = 7676.571599999999
$ws.Range("L122").Value = 17884.9995
$ws.Range("M122").Value = -5226.571599999999
$ws.Range("N122").Value = -22784.9995
$ws.Range("H126").Value = 3579.1155
$ws.Range("I126").Value = 2294.2354
$ws.Range("J126").Value = 6006.1113
$ws.Range("K126").Value = 6882.706200000001
$ws.Range("L126").Value = 18018.3339
$ws.Range("M126").Value = -4412.706200000001
$ws.Range("N126").Value = -22958.3339
$ws.Range("H127").Value = 68888.78
$ws.Range("J127").Value = 68888.78
$ws.Range("L127").Value = 68888.78
$ws.Range("N127").Value = -78808.78
$ws.Range("H132").Value = 3609.1794
$ws.Range("I132").Value = 2792.4614
$ws.Range("K132").Value = 8377.3842
$ws.Range("M132").Value = -5847.3842
$ws.Range("H135").Value = 80000
$ws.Range("J135").Value = 80000
$ws.Range("L135").Value = 80000
$ws.Range("N135").Value = -90140
$ws.Range("H136").Value = 4966.0713
$ws.Range("I136").Value = 3164.7646
$ws.Range("J136").Value = 7749.909
$ws.Range("K136").Value = 9494.293799999999
$ws.Range("L136").Value = 23249.727
$ws.Range("M136").Value = -6944.293799999999
$ws.Range("N136").Value = -28349.727
$ws.Range("H137").Value = 67894.37
$ws.Range("J137").Value = 69721.836
$ws.Range("L137").Value = 69721.836
$ws.Range("N137").Value = -79921.836

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 38495
$ws.Range("J54").Value = 38495
$ws.Range("L54").Value = 38495
$ws.Range("N54").Value = -39535
$ws.Range("H68").Value = 29499.5
$ws.Range("J68").Value = 29499.5
$ws.Range("L68").Value = 29499.5
$ws.Range("N68").Value = -31121.5
$ws.Range("H71").Value = 29499.5
$ws.Range("J71").Value = 29499.5
$ws.Range("L71").Value = 88498.5
$ws.Range("N71").Value = -96610.5
$ws.Range("H81").Value = 7558.4707
$ws.Range("I81").Value = 12643.777
$ws.Range("J81").Value = 1837.5
$ws.Range("K81").Value = 25287.554
$ws.Range("L81").Value = 3675
$ws.Range("M81").Value = -24226.554
$ws.Range("N81").Value = -5797
$ws.Range("H84").Value = 7558.4707
$ws.Range("I84").Value = 12643.777
$ws.Range("J84").Value = 1837.5
$ws.Range("K84").Value = 126437.77
$ws.Range("L84").Value = 18375
$ws.Range("M84").Value = -121133.77
$ws.Range("N84").Value = -28983
$ws.Range("H100").Value = 712.56525
$ws.Range("I100").Value = 599.2222
$ws.Range("K100").Value = 1198.4444
$ws.Range("M100").Value = -657.4444000000001
$ws.Range("H113").Value = 349.9
$ws.Range("I113").Value = 272.7143
$ws.Range("J113").Value = 530
$ws.Range("K113").Value = 818.1428999999999
$ws.Range("L113").Value = 1590
$ws.Range("M113").Value = 1351.8571
$ws.Range("N113").Value = -5930
$ws.Range("H132").Value = 1974.1277
$ws.Range("I132").Value = 1229.7949
$ws.Range("J132").Value = 5602.75
$ws.Range("K132").Value = 3689.384700000001
$ws.Range("L132").Value = 16808.25
$ws.Range("M132").Value = -1159.384700000001
$ws.Range("N132").Value = -21868.25
$ws.Range("H136").Value = 2140.7874
$ws.Range("I136").Value = 1181.0834
$ws.Range("J136").Value = 5281.636
$ws.Range("K136").Value = 3543.2502
$ws.Range("L136").Value = 15844.908
$ws.Range("M136").Value = -993.2501999999999
$ws.Range("N136").Value = -20944.908
